$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert the three new <comment> / c_038v_01 / </comment> runs right after
#    "leur lb" and before " de vente &" (the run boundary straddling the
#    commentReference in the paragraph about "de vente").
# ---------------------------------------------------------------------------
$anchor = $d.Content
$found = $anchor.Find.Execute("leur lb de vente", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Host "ERROR: anchor text not found"
}

# "leur lb" is 7 characters -> position right after the comment-referenced "b"
$insStart = $anchor.Start + 7

$text1 = "<comment>"
$text2 = "c_038v_01"
$text3 = "</comment>"

$ins = $d.Range($insStart, $insStart)
$ins.InsertAfter($text1 + $text2 + $text3)

$len1 = $text1.Length
$len2 = $text2.Length
$len3 = $text3.Length

$r1 = $d.Range($insStart, $insStart + $len1)
$r2 = $d.Range($insStart + $len1, $insStart + $len1 + $len2)
$r3 = $d.Range($insStart + $len1 + $len2, $insStart + $len1 + $len2 + $len3)

# Run 1: "<comment>" - Courier New, blue, 9pt
$r1.Font.Name = "Courier New"
$r1.Font.Color = 16711680
$r1.Font.Size = 9

# Run 2: "c_038v_01" - 8.5pt, white highlight
$r2.Font.Size = 8.5

# Run 3: "</comment>" - Courier New, blue, 9pt, white highlight
$r3.Font.Name = "Courier New"
$r3.Font.Color = 16711680
$r3.Font.Size = 9

# Apply the white highlight via a range-scoped Find/Replace (the only path in
# this host that actually records w:highlight on the run).
$r2h = $d.Range($insStart + $len1, $insStart + $len1 + $len2)
$r2h.Find.ClearFormatting()
$r2h.Find.Replacement.ClearFormatting()
$r2h.Find.Replacement.Highlight = 8
$r2h.Find.Text = $text2
$r2h.Find.Replacement.Text = "^&"
$r2h.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2, $true, $null, $null, $null)

$r3h = $d.Range($insStart + $len1 + $len2, $insStart + $len1 + $len2 + $len3)
$r3h.Find.ClearFormatting()
$r3h.Find.Replacement.ClearFormatting()
$r3h.Find.Replacement.Highlight = 8
$r3h.Find.Text = $text3
$r3h.Find.Replacement.Text = "^&"
$r3h.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2, $true, $null, $null, $null)

# ---------------------------------------------------------------------------
# 2) sectPr: add a footer distance of 720 twips (0.5in = 36pt) to pgMar.
# ---------------------------------------------------------------------------
$d.PageSetup.FooterDistance = 36

Write-Host "edit complete"
